# Keldysh energy window: add a new "Sheet2" with the Hamiltonian/eigenvector
# matrix-size breakdown (mirrors the layout already used on Sheet1), and make
# Sheet2 the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Re-apply the numeric format on Sheet1!C8:E8 so that it keeps the "0.0"
# display format used by the rest of the Size-in-GB row.
$ws1.Range("C8:E8").NumberFormat = "0.0"

# Add the new worksheet right after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# --- Symmetric block -----------------------------------------------------
$ws2.Range("A1").Value = "Ham/Eigvecs"
$ws2.Range("B1").Value = 4
$ws2.Range("C1").Value = 2

$ws2.Range("D2").Value = "G Sym"
$ws2.Range("E2").Value = 2
$ws2.Range("F2").Value = 1

$ws2.Range("G3").Value = "Z"
$ws2.Range("H3").Value = 4
$ws2.Range("I3").Value = 1

$ws2.Range("G4").Value = "A"
$ws2.Range("H4").Value = 2
$ws2.Range("I4").Value = 1

$ws2.Range("G6").Value = "C_ab"
$ws2.Range("H6").Value = 1
$ws2.Range("I6").Value = 2

$ws2.Range("G7").Value = "C_gg"
$ws2.Range("H7").Value = 1
$ws2.Range("I7").Value = 2

$ws2.Range("G8").Value = "output"
$ws2.Range("H8").Value = 1
$ws2.Range("I8").Value = 2

# --- Asymmetric block ------------------------------------------------------
$ws2.Range("D10").Value = "G Asym"
$ws2.Range("E10").Value = 2
$ws2.Range("F10").Value = 2

$ws2.Range("G11").Value = "Z"
$ws2.Range("H11").Value = 4
$ws2.Range("I11").Value = 2

$ws2.Range("G12").Value = "A"
$ws2.Range("H12").Value = 2
$ws2.Range("I12").Value = 2

$ws2.Range("G14").Value = "C_ab"
$ws2.Range("H14").Value = 1
$ws2.Range("I14").Value = 2
$ws2.Range("H14:I14").Locked = $true

$ws2.Range("G15").Value = "C_gg"
$ws2.Range("H15").Value = 1
$ws2.Range("I15").Value = 2
$ws2.Range("H15:I15").Locked = $true

$ws2.Range("G16").Value = "output"
$ws2.Range("H16").Value = 1
$ws2.Range("I16").Value = 2
$ws2.Range("H16:I16").Locked = $true

# --- Summary table ---------------------------------------------------------
$ws2.Range("A20").Value = "H"
$ws2.Range("B20").Value = "G"
$ws2.Range("C20").Value = "Z"
$ws2.Range("D20").Value = "A"
$ws2.Range("E20").Value = "C_ab"
$ws2.Range("F20").Value = "C_gg"
$ws2.Range("G20").Value = "output"

$ws2.Range("A21").Value = 8
$ws2.Range("B21").Value = 2
$ws2.Range("C21").Value = 4
$ws2.Range("D21").Value = 2
$ws2.Range("H21").Formula = "=SUM(A21:G21)"

$ws2.Range("A22").Value = 8
$ws2.Range("B22").Value = 2
$ws2.Range("E22").Value = 2
$ws2.Range("F22").Value = 2
$ws2.Range("G22").Value = 2
$ws2.Range("H22").Formula = "=SUM(A22:G22)"

$ws2.Range("A23").Value = 8
$ws2.Range("B23").Value = 4
$ws2.Range("C23").Value = 8
$ws2.Range("D23").Value = 4
$ws2.Range("H23").Formula = "=SUM(A23:G23)"

$ws2.Range("A24").Value = 8
$ws2.Range("B24").Value = 4
$ws2.Range("E24").Value = 2
$ws2.Range("F24").Value = 2
$ws2.Range("G24").Value = 2
$ws2.Range("H24").Formula = "=SUM(A24:G24)"

# Sheet2 is the active tab.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 100
